# Reroute test data update
# Regenerates the random "prodXXXX" ProductName values for the three rows
# (Sku rows 2, 3 and 5 on the "Input" sheet) whose names are being rotated
# in this pass, and re-stamps their cell format (white solid fill + thin
# top/bottom border) so the row keeps the same look it always had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-ProductNameCell {
    param(
        [string]$CellRef,
        [string]$NewName
    )

    $cell = $ws.Range($CellRef)

    # New random product name for this refresh.
    $cell.Value = $NewName

    # Re-apply the row's standard look: solid white fill, thin top & bottom
    # border -- matches every other ProductName cell in the sheet.
    $cell.Interior.ColorIndex = 2

    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
}

Set-ProductNameCell "B3" "prodhQod"
Set-ProductNameCell "B5" "prodXHSv"
Set-ProductNameCell "B2" "prodaCWk"
